# =====================================================================
# PlayerPerformance_4420.xlsx restructuring
#
#  1. Insert a new 'Player Info' sheet at the front (ID / NAME /
#     BATTING_HAND / BOWL_STYLE) describing the player.
#  2. 'ODI Batting': rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#     full howstat.com URL with the bare numeric match code. Also drop the
#     placeholder empty INNING_NUMBER cells on 'did not bat' rows.
#  3. 'ODI Bowling': rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#     full howstat.com URL with the bare numeric match code.
#  4. Append a new 'ODI Batting Extra' sheet at the end with additional
#     per-match batting detail (position, boundaries, % of total, MoM).
# =====================================================================

$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# ---- 1. 'Player Info' sheet (new first sheet) ----
$battingSheetBefore = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($battingSheetBefore)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID","NAME","BATTING_HAND","BOWL_STYLE")
for ($col = 1; $col -le $playerInfoHeaders.Length; $col++) {
    $playerInfo.Cells.Item(1, $col).Value = $playerInfoHeaders[$col - 1]
}
Set-HeaderStyle($playerInfo.Range("A1:D1"))

$playerInfo.Cells.Item(2, 1).Value = "'4420"
$playerInfo.Cells.Item(2, 2).Value = "David J Willey"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Left Arm Fast Medium"

# ---- 2. 'ODI Batting' sheet: MATCH_CARD_LINK -> MATCH_CODE ----
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingMatchCodes = @("3800","3806","3807","3809","3833","3834","3858","3859","3861","3863","3885","3887","3889","3904","3906","3908","3910","3911","3930","3932","3946","3948","3974","3976","3978","4020","4021","4030","4068","4070","4072","4125","4137","4138","4165","4166","4167","4168","4169","4171","4173","4175","4284","4292","4294","4300","4426","4427","4428","4469","4470","4471","4598","4599","4602","4609","4613","4618","4620","4622","4660","4663","4666","4698")
$row = 2
foreach ($code in $battingMatchCodes) {
    $battingSheet.Cells.Item($row, 4).Value = "'" + $code
    $row = $row + 1
}

# Rows where the player did not bat leave INNING_NUMBER (col B) blank;
# the source cell is removed entirely rather than kept as an empty cell.
$battingEmptyInningRows = @(2,4,9,10,11,13,16,17,18,19,26,27,30,31,35,40,43,45,46,48,51,52,53,54,55,56,61)
foreach ($r in $battingEmptyInningRows) {
    $battingSheet.Cells.Item($r, 2).Value = ""
}

# ---- 3. 'ODI Bowling' sheet: MATCH_CARD_LINK -> MATCH_CODE ----
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingMatchCodes = @("3800","3806","3807","3809","3833","3834","3858","3859","3861","3863","3885","3887","3889","3904","3906","3908","3910","3911","3930","3932","3946","3948","3974","3976","3978","4020","4021","4030","4068","4072","4125","4137","4138","4165","4166","4167","4168","4169","4171","4173","4175","4284","4292","4294","4300","4426","4427","4428","4469","4470","4471","4598","4599","4602","4609","4613","4618","4620","4622","4660","4663","4666","4698")
$row = 2
foreach ($code in $bowlingMatchCodes) {
    $bowlingSheet.Cells.Item($row, 2).Value = "'" + $code
    $row = $row + 1
}

# ---- 4. 'ODI Batting Extra' sheet (new last sheet) ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtraHeaders = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($col = 1; $col -le $battingExtraHeaders.Length; $col++) {
    $battingExtra.Cells.Item(1, $col).Value = $battingExtraHeaders[$col - 1]
}
Set-HeaderStyle($battingExtra.Range("A1:F1"))

# MATCH_CODE (text), BATTING_POSITION (number or blank), NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL (text or blank), MAN_OF_MATCH (text)
$battingExtra.Cells.Item(2, 1).Value = "'4294"
$battingExtra.Cells.Item(2, 2).Value = ""
$battingExtra.Cells.Item(2, 3).Value = ""
$battingExtra.Cells.Item(2, 4).Value = ""
$battingExtra.Cells.Item(2, 5).Value = ""
$battingExtra.Cells.Item(2, 6).Value = "NO"
$battingExtra.Cells.Item(3, 1).Value = "'4300"
$battingExtra.Cells.Item(3, 2).Value = 9
$battingExtra.Cells.Item(3, 3).Value = "'0"
$battingExtra.Cells.Item(3, 4).Value = "'1"
$battingExtra.Cells.Item(3, 5).Value = "'3.99%"
$battingExtra.Cells.Item(3, 6).Value = "NO"
$battingExtra.Cells.Item(4, 1).Value = "'4426"
$battingExtra.Cells.Item(4, 2).Value = ""
$battingExtra.Cells.Item(4, 3).Value = ""
$battingExtra.Cells.Item(4, 4).Value = ""
$battingExtra.Cells.Item(4, 5).Value = ""
$battingExtra.Cells.Item(4, 6).Value = "NO"
$battingExtra.Cells.Item(5, 1).Value = "'4427"
$battingExtra.Cells.Item(5, 2).Value = 8
$battingExtra.Cells.Item(5, 3).Value = "'5"
$battingExtra.Cells.Item(5, 4).Value = "'2"
$battingExtra.Cells.Item(5, 5).Value = "'21.76%"
$battingExtra.Cells.Item(5, 6).Value = "NO"
$battingExtra.Cells.Item(6, 1).Value = "'4428"
$battingExtra.Cells.Item(6, 2).Value = 8
$battingExtra.Cells.Item(6, 3).Value = "'3"
$battingExtra.Cells.Item(6, 4).Value = "'3"
$battingExtra.Cells.Item(6, 5).Value = "'15.55%"
$battingExtra.Cells.Item(6, 6).Value = "NO"
$battingExtra.Cells.Item(7, 1).Value = "'4469"
$battingExtra.Cells.Item(7, 2).Value = 9
$battingExtra.Cells.Item(7, 3).Value = ""
$battingExtra.Cells.Item(7, 4).Value = ""
$battingExtra.Cells.Item(7, 5).Value = ""
$battingExtra.Cells.Item(7, 6).Value = "NO"
$battingExtra.Cells.Item(8, 1).Value = "'4470"
$battingExtra.Cells.Item(8, 2).Value = 9
$battingExtra.Cells.Item(8, 3).Value = ""
$battingExtra.Cells.Item(8, 4).Value = ""
$battingExtra.Cells.Item(8, 5).Value = ""
$battingExtra.Cells.Item(8, 6).Value = "NO"
$battingExtra.Cells.Item(9, 1).Value = "'4471"
$battingExtra.Cells.Item(9, 2).Value = ""
$battingExtra.Cells.Item(9, 3).Value = ""
$battingExtra.Cells.Item(9, 4).Value = ""
$battingExtra.Cells.Item(9, 5).Value = ""
$battingExtra.Cells.Item(9, 6).Value = "NO"
$battingExtra.Cells.Item(10, 1).Value = "'4598"
$battingExtra.Cells.Item(10, 2).Value = ""
$battingExtra.Cells.Item(10, 3).Value = ""
$battingExtra.Cells.Item(10, 4).Value = ""
$battingExtra.Cells.Item(10, 5).Value = ""
$battingExtra.Cells.Item(10, 6).Value = "NO"
$battingExtra.Cells.Item(11, 1).Value = "'4599"
$battingExtra.Cells.Item(11, 2).Value = ""
$battingExtra.Cells.Item(11, 3).Value = ""
$battingExtra.Cells.Item(11, 4).Value = ""
$battingExtra.Cells.Item(11, 5).Value = ""
$battingExtra.Cells.Item(11, 6).Value = "NO"
$battingExtra.Cells.Item(12, 1).Value = "'4602"
$battingExtra.Cells.Item(12, 2).Value = 8
$battingExtra.Cells.Item(12, 3).Value = ""
$battingExtra.Cells.Item(12, 4).Value = ""
$battingExtra.Cells.Item(12, 5).Value = ""
$battingExtra.Cells.Item(12, 6).Value = "NO"
$battingExtra.Cells.Item(13, 1).Value = "'4609"
$battingExtra.Cells.Item(13, 2).Value = 8
$battingExtra.Cells.Item(13, 3).Value = "'3"
$battingExtra.Cells.Item(13, 4).Value = "'0"
$battingExtra.Cells.Item(13, 5).Value = "'19.09%"
$battingExtra.Cells.Item(13, 6).Value = "NO"
$battingExtra.Cells.Item(14, 1).Value = "'4613"
$battingExtra.Cells.Item(14, 2).Value = 8
$battingExtra.Cells.Item(14, 3).Value = "'2"
$battingExtra.Cells.Item(14, 4).Value = "'2"
$battingExtra.Cells.Item(14, 5).Value = "'16.67%"
$battingExtra.Cells.Item(14, 6).Value = "NO"
$battingExtra.Cells.Item(15, 1).Value = "'4618"
$battingExtra.Cells.Item(15, 2).Value = 8
$battingExtra.Cells.Item(15, 3).Value = "'1"
$battingExtra.Cells.Item(15, 4).Value = "'1"
$battingExtra.Cells.Item(15, 5).Value = "'6.95%"
$battingExtra.Cells.Item(15, 6).Value = "NO"
$battingExtra.Cells.Item(16, 1).Value = "'4620"
$battingExtra.Cells.Item(16, 2).Value = 9
$battingExtra.Cells.Item(16, 3).Value = "'1"
$battingExtra.Cells.Item(16, 4).Value = "'1"
$battingExtra.Cells.Item(16, 5).Value = "'10.45%"
$battingExtra.Cells.Item(16, 6).Value = "NO"
$battingExtra.Cells.Item(17, 1).Value = "'4622"
$battingExtra.Cells.Item(17, 2).Value = ""
$battingExtra.Cells.Item(17, 3).Value = ""
$battingExtra.Cells.Item(17, 4).Value = ""
$battingExtra.Cells.Item(17, 5).Value = ""
$battingExtra.Cells.Item(17, 6).Value = "NO"
$battingExtra.Cells.Item(18, 1).Value = "'4660"
$battingExtra.Cells.Item(18, 2).Value = ""
$battingExtra.Cells.Item(18, 3).Value = ""
$battingExtra.Cells.Item(18, 4).Value = ""
$battingExtra.Cells.Item(18, 5).Value = ""
$battingExtra.Cells.Item(18, 6).Value = "NO"
$battingExtra.Cells.Item(19, 1).Value = "'4663"
$battingExtra.Cells.Item(19, 2).Value = ""
$battingExtra.Cells.Item(19, 3).Value = ""
$battingExtra.Cells.Item(19, 4).Value = ""
$battingExtra.Cells.Item(19, 5).Value = ""
$battingExtra.Cells.Item(19, 6).Value = "NO"
$battingExtra.Cells.Item(20, 1).Value = "'4666"
$battingExtra.Cells.Item(20, 2).Value = ""
$battingExtra.Cells.Item(20, 3).Value = ""
$battingExtra.Cells.Item(20, 4).Value = ""
$battingExtra.Cells.Item(20, 5).Value = ""
$battingExtra.Cells.Item(20, 6).Value = "NO"
$battingExtra.Cells.Item(21, 1).Value = "'4698"
$battingExtra.Cells.Item(21, 2).Value = 8
$battingExtra.Cells.Item(21, 3).Value = "'2"
$battingExtra.Cells.Item(21, 4).Value = "'0"
$battingExtra.Cells.Item(21, 5).Value = "'2.95%"
$battingExtra.Cells.Item(21, 6).Value = "NO"

Write-Output "Workbook restructured: Player Info / ODI Batting / ODI Bowling / ODI Batting Extra"
